# Update attendance/view counts (column F) on the "展览" and "全部类型"
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 615
$ws1.Range("F9").Value  = 4476
$ws1.Range("F10").Value = 1803
$ws1.Range("F12").Value = 122
$ws1.Range("F13").Value = 3017
$ws1.Range("F17").Value = 567
$ws1.Range("F18").Value = 487
$ws1.Range("F35").Value = 3257
$ws1.Range("F38").Value = 196
$ws1.Range("F40").Value = 1232

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 615
$ws4.Range("F9").Value  = 4476
$ws4.Range("F10").Value = 1803
$ws4.Range("F12").Value = 122
$ws4.Range("F13").Value = 3017
$ws4.Range("F17").Value = 567
$ws4.Range("F18").Value = 487
$ws4.Range("F36").Value = 3257
$ws4.Range("F40").Value = 196
$ws4.Range("F42").Value = 1232
